$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FirstName/LastName values in row 4 ("Z" / "Man" -> "Zman" / "Zach")
$ws.Range("A4").Value = "Zman"
$ws.Range("B4").Value = "Zach"

# Move the active selection from C4 to B4
$ws.Range("B4").Select()
